$wb = $excel.ActiveWorkbook
$wsAVL = $wb.Worksheets.Item("AVLRaPTC")

# Update C7 on the AVLRaPTC sheet: it now references C3 instead of a literal 0
$wsAVL.Range("C7").Formula = "=C3"

# The workbook was last saved with AVLRaPTC as the active sheet, with C8 selected
$wsAVL.Activate() | Out-Null
$wsAVL.Range("C8").Select() | Out-Null
